$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting (style) from an existing date cell (B5) onto B23
# so the new date reuses the existing numFmt/style instead of creating a new one.
$ws.Range("B5").Copy()
$ws.Range("B23").PasteSpecial(-4122)

# Fill in the new time-tracking entry in row 23
$ws.Range("A23").Value = 1.5
$ws.Range("B23").Value = 43549
$ws.Range("C23").Value = "16:30-18:00"
$ws.Range("D23").Value = "Prästation - Übung 3"

# Update the active selection to match the authored state
$ws.Range("I8").Select() | Out-Null
